$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update modified timestamp in B20
$ws.Range("B20").Value = "2022-06-20T06:20:30+00:00"

# Fill in dates/creator for row 23 (vocab:1000 - reaction time)
$ws.Range("T23").NumberFormat = "@"
$ws.Range("T23").Value = "2022-06-20"
$ws.Range("T23").Style = "Normal"

$ws.Range("U23").NumberFormat = "@"
$ws.Range("U23").Value = "2022-06-20"
$ws.Range("U23").Style = "Normal"

$ws.Range("V23").Value = "0000-0003-2195-3997"

# Fill in dates/creator for row 24 (vocab:1001 - reaction time at prospective memory task)
$ws.Range("T24").NumberFormat = "@"
$ws.Range("T24").Value = "2022-06-20"
$ws.Range("T24").Style = "Normal"

$ws.Range("U24").NumberFormat = "@"
$ws.Range("U24").Value = "2022-06-20"
$ws.Range("U24").Style = "Normal"

$ws.Range("V24").Value = "0000-0003-2195-3997"
